$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J width (target stored width 36.44140625 chars; engine rounds
# ColumnWidth to the nearest 1/6 character, so 35.6 is the closest input
# that reproduces the nearest representable stored width)
$ws.Columns("J").ColumnWidth = 35.6

# Header cell J5: same formatting as the other header cells (copy format from B5)
$ws.Range("B5").Copy() | Out-Null
$ws.Range("J5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("J5").Value = "No standardization`nshuffle when use kfold (seed = 1)"

# Data cells J6:J10: same formatting as column B data cells
$ws.Range("B6").Copy() | Out-Null
$ws.Range("J6:J10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("J6").Value = 4.8006113476114498
$ws.Range("J7").Value = 4.8002483703501797
$ws.Range("J8").Value = 4.8127441256846701
$ws.Range("J9").Value = 4.8634390971983104
$ws.Range("J10").Value = 4.9566232289732701

# Update sheet view: scroll position (best effort) and active selection cell
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 5
$ws.Range("J18").Select()
